$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the current used range extent before we insert anything
$lastRow = $ws.UsedRange.Rows.Count

# Insert a new column before column A, shifting existing data (A:G) to (B:H)
$ws.Range("A1").EntireColumn.Insert()

# New header cell — matches the bold/centered style used by the other header cells
$ws.Range("A1").Value = "select"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").HorizontalAlignment = -4108

# Fill all data rows (2..lastRow) in the new column A with "Y" (select everything)
$ws.Range("A2:A" + $lastRow).Value = "Y"

# Select the full updated range
$ws.Range("A1:H" + $lastRow).Select()
